$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: snapshot current B:AD values for every affected row (before any writes)
$snapshot = @{}
$snapshot[25] = $ws.Range("B25:AD25").Value()
$snapshot[26] = $ws.Range("B26:AD26").Value()
$snapshot[64] = $ws.Range("B64:AD64").Value()
$snapshot[65] = $ws.Range("B65:AD65").Value()
$snapshot[66] = $ws.Range("B66:AD66").Value()
$snapshot[67] = $ws.Range("B67:AD67").Value()
$snapshot[101] = $ws.Range("B101:AD101").Value()
$snapshot[102] = $ws.Range("B102:AD102").Value()
$snapshot[119] = $ws.Range("B119:AD119").Value()
$snapshot[120] = $ws.Range("B120:AD120").Value()
$snapshot[121] = $ws.Range("B121:AD121").Value()
$snapshot[130] = $ws.Range("B130:AD130").Value()
$snapshot[131] = $ws.Range("B131:AD131").Value()
$snapshot[132] = $ws.Range("B132:AD132").Value()
$snapshot[133] = $ws.Range("B133:AD133").Value()
$snapshot[141] = $ws.Range("B141:AD141").Value()
$snapshot[142] = $ws.Range("B142:AD142").Value()
$snapshot[145] = $ws.Range("B145:AD145").Value()
$snapshot[146] = $ws.Range("B146:AD146").Value()
$snapshot[147] = $ws.Range("B147:AD147").Value()
$snapshot[215] = $ws.Range("B215:AD215").Value()
$snapshot[216] = $ws.Range("B216:AD216").Value()
$snapshot[217] = $ws.Range("B217:AD217").Value()
$snapshot[240] = $ws.Range("B240:AD240").Value()
$snapshot[241] = $ws.Range("B241:AD241").Value()
$snapshot[254] = $ws.Range("B254:AD254").Value()
$snapshot[255] = $ws.Range("B255:AD255").Value()
$snapshot[259] = $ws.Range("B259:AD259").Value()
$snapshot[261] = $ws.Range("B261:AD261").Value()

# Step 2: write back according to the row-swap mapping derived from the commit diff
$ws.Range("B25:AD25").Value = $snapshot[26]
$ws.Range("B26:AD26").Value = $snapshot[25]
$ws.Range("B64:AD64").Value = $snapshot[67]
$ws.Range("B65:AD65").Value = $snapshot[64]
$ws.Range("B66:AD66").Value = $snapshot[65]
$ws.Range("B67:AD67").Value = $snapshot[66]
$ws.Range("B101:AD101").Value = $snapshot[102]
$ws.Range("B102:AD102").Value = $snapshot[101]
$ws.Range("B119:AD119").Value = $snapshot[121]
$ws.Range("B120:AD120").Value = $snapshot[119]
$ws.Range("B121:AD121").Value = $snapshot[120]
$ws.Range("B130:AD130").Value = $snapshot[131]
$ws.Range("B131:AD131").Value = $snapshot[130]
$ws.Range("B132:AD132").Value = $snapshot[133]
$ws.Range("B133:AD133").Value = $snapshot[132]
$ws.Range("B141:AD141").Value = $snapshot[142]
$ws.Range("B142:AD142").Value = $snapshot[141]
$ws.Range("B145:AD145").Value = $snapshot[147]
$ws.Range("B146:AD146").Value = $snapshot[145]
$ws.Range("B147:AD147").Value = $snapshot[146]
$ws.Range("B215:AD215").Value = $snapshot[216]
$ws.Range("B216:AD216").Value = $snapshot[217]
$ws.Range("B217:AD217").Value = $snapshot[215]
$ws.Range("B240:AD240").Value = $snapshot[241]
$ws.Range("B241:AD241").Value = $snapshot[240]
$ws.Range("B254:AD254").Value = $snapshot[255]
$ws.Range("B255:AD255").Value = $snapshot[254]
$ws.Range("B259:AD259").Value = $snapshot[261]
$ws.Range("B261:AD261").Value = $snapshot[259]

Write-Host "Row data permutation applied successfully."
